$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 58
$ws.Range("I2").Value = 129
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 167
$ws.Range("N2").Value = 92
$ws.Range("P2").Value = 7
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 73
$ws.Range("T2").Value = 106
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 883
$ws.Range("X2").Value = 888
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 4
